$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# Row 53 (ALC) - hunk 0
$ws_ALC.Range("H53").Value = 303.26315
$ws_ALC.Range("I53").Value = 184.28572
$ws_ALC.Range("J53").Value = 636.4
$ws_ALC.Range("K53").Value = 184.28572
$ws_ALC.Range("L53").Value = 636.4
$ws_ALC.Range("M53").Value = 452.71428
$ws_ALC.Range("N53").Value = -1910.4

# Row 62 (ALC) - hunk 1
$ws_ALC.Range("H62").Value = 306776.47
$ws_ALC.Range("I62").Value = 479548
$ws_ALC.Range("K62").Value = 479548
$ws_ALC.Range("M62").Value = -478924

# Row 65 (ALC) - hunk 2
$ws_ALC.Range("H65").Value = 306776.47
$ws_ALC.Range("I65").Value = 479548
$ws_ALC.Range("K65").Value = 2397740
$ws_ALC.Range("M65").Value = -2394620

# Row 86 (ALC) - hunk 3
$ws_ALC.Range("H86").Value = 125005200
$ws_ALC.Range("I86").Value = 4650.75
$ws_ALC.Range("J86").Value = 250005740
$ws_ALC.Range("K86").Value = 4650.75
$ws_ALC.Range("L86").Value = 250005740
$ws_ALC.Range("M86").Value = -3527.75
$ws_ALC.Range("N86").Value = -250007986

# Row 89 (ALC) - hunk 4
$ws_ALC.Range("H89").Value = 125005200
$ws_ALC.Range("I89").Value = 4650.75
$ws_ALC.Range("J89").Value = 250005740
$ws_ALC.Range("K89").Value = 23253.75
$ws_ALC.Range("L89").Value = 1250028700
$ws_ALC.Range("M89").Value = -17637.75
$ws_ALC.Range("N89").Value = -1250039932

# Row 98 (ALC) - hunk 5
$ws_ALC.Range("H98").Value = 1587.4445
$ws_ALC.Range("I98").Value = 1682.2
$ws_ALC.Range("J98").Value = 403
$ws_ALC.Range("K98").Value = 1682.2
$ws_ALC.Range("L98").Value = 403
$ws_ALC.Range("M98").Value = -184.2
$ws_ALC.Range("N98").Value = -3399

# Row 116 (ALC) - hunk 6
$ws_ALC.Range("H116").Value = 4046.4348
$ws_ALC.Range("J116").Value = 4333
$ws_ALC.Range("L116").Value = 4333
$ws_ALC.Range("N116").Value = -11217

# Row 122 (ALC) - hunk 7
$ws_ALC.Range("H122").Value = 1587.4445
$ws_ALC.Range("I122").Value = 1682.2
$ws_ALC.Range("J122").Value = 403
$ws_ALC.Range("K122").Value = 5046.6
$ws_ALC.Range("L122").Value = 1209
$ws_ALC.Range("M122").Value = -2596.6
$ws_ALC.Range("N122").Value = -6109

# Row 137 (ALC) - hunk 8
$ws_ALC.Range("H137").Value = 1040.5555
$ws_ALC.Range("I137").Value = 987
$ws_ALC.Range("J137").Value = 1179.8
$ws_ALC.Range("K137").Value = 2961
$ws_ALC.Range("L137").Value = 3539.4
$ws_ALC.Range("M137").Value = -411
$ws_ALC.Range("N137").Value = -8639.4

# Row 32 (ARM) - hunk 9
$ws_ARM.Range("H32").Value = 2249.14
$ws_ARM.Range("I32").Value = 2249.14
$ws_ARM.Range("K32").Value = 2249.14
$ws_ARM.Range("M32").Value = -1962.14

# Row 63 (ARM) - hunk 10
$ws_ARM.Range("H63").Value = 2052.4375
$ws_ARM.Range("I63").Value = 2052.4375
$ws_ARM.Range("J63").Value = 0
$ws_ARM.Range("K63").Value = 2052.4375
$ws_ARM.Range("L63").Value = 0
$ws_ARM.Range("M63").Value = -1366.4375
$ws_ARM.Range("N63").ClearContents()

# Row 66 (ARM) - hunk 11
$ws_ARM.Range("H66").Value = 2052.4375
$ws_ARM.Range("I66").Value = 2052.4375
$ws_ARM.Range("J66").Value = 0
$ws_ARM.Range("K66").Value = 10262.1875
$ws_ARM.Range("L66").Value = 0
$ws_ARM.Range("M66").Value = -6830.1875
$ws_ARM.Range("N66").ClearContents()

# Row 122 (ARM) - hunk 12
$ws_ARM.Range("H122").Value = 1620.2
$ws_ARM.Range("I122").Value = 1300.5
$ws_ARM.Range("K122").Value = 3901.5
$ws_ARM.Range("M122").Value = -1451.5

# Row 20 (BSM) - hunk 13
$ws_BSM.Range("H20").Value = 27799932
$ws_BSM.Range("I20").Value = 29240.521
$ws_BSM.Range("J20").Value = 76932696
$ws_BSM.Range("K20").Value = 29240.521
$ws_BSM.Range("L20").Value = 76932696
$ws_BSM.Range("M20").Value = -28993.521
$ws_BSM.Range("N20").Value = -76933190

# Row 105 (BSM) - hunk 14
$ws_BSM.Range("H105").Value = 3825.3
$ws_BSM.Range("I105").Value = 3853.3794
$ws_BSM.Range("K105").Value = 3853.3794
$ws_BSM.Range("M105").Value = -2106.3794

# Row 31 (CRP) - hunk 15
$ws_CRP.Range("H31").Value = 76581.57000000001
$ws_CRP.Range("I31").Value = 6512
$ws_CRP.Range("J31").Value = 202706.8
$ws_CRP.Range("K31").Value = 6512
$ws_CRP.Range("L31").Value = 202706.8
$ws_CRP.Range("M31").Value = -6217
$ws_CRP.Range("N31").Value = -203296.8

# Row 34 (CRP) - hunk 16
$ws_CRP.Range("H34").Value = 76581.57000000001
$ws_CRP.Range("I34").Value = 6512
$ws_CRP.Range("J34").Value = 202706.8
$ws_CRP.Range("K34").Value = 6512
$ws_CRP.Range("L34").Value = 202706.8
$ws_CRP.Range("M34").Value = -6310
$ws_CRP.Range("N34").Value = -203110.8

# Row 99 (CRP) - hunk 17
$ws_CRP.Range("H99").Value = 3426.4866
$ws_CRP.Range("I99").Value = 2704.1365
$ws_CRP.Range("J99").Value = 4485.933
$ws_CRP.Range("K99").Value = 2704.1365
$ws_CRP.Range("L99").Value = 4485.933
$ws_CRP.Range("M99").Value = -1206.1365
$ws_CRP.Range("N99").Value = -7481.933

# Row 105 (CRP) - hunk 18
$ws_CRP.Range("H105").Value = 1811.1111
$ws_CRP.Range("I105").Value = 1811.1111
$ws_CRP.Range("K105").Value = 1811.1111
$ws_CRP.Range("M105").Value = -64.11110000000008

# Row 126 (CRP) - hunk 19
$ws_CRP.Range("H126").Value = 3426.4866
$ws_CRP.Range("I126").Value = 2704.1365
$ws_CRP.Range("J126").Value = 4485.933
$ws_CRP.Range("K126").Value = 8112.4095
$ws_CRP.Range("L126").Value = 13457.799
$ws_CRP.Range("M126").Value = -5642.4095
$ws_CRP.Range("N126").Value = -18397.799

# Row 70 (GSM) - hunk 20
$ws_GSM.Range("H70").Value = 4539.2
$ws_GSM.Range("I70").Value = 4084.111
$ws_GSM.Range("J70").Value = 4911.5454
$ws_GSM.Range("K70").Value = 4084.111
$ws_GSM.Range("L70").Value = 4911.5454
$ws_GSM.Range("M70").Value = -3814.111
$ws_GSM.Range("N70").Value = -5451.5454

# Row 73 (GSM) - hunk 21
$ws_GSM.Range("H73").Value = 4539.2
$ws_GSM.Range("I73").Value = 4084.111
$ws_GSM.Range("J73").Value = 4911.5454
$ws_GSM.Range("K73").Value = 4084.111
$ws_GSM.Range("L73").Value = 4911.5454
$ws_GSM.Range("M73").Value = -3148.111
$ws_GSM.Range("N73").Value = -6783.5454

# Row 122 (GSM) - hunk 22
$ws_GSM.Range("H122").Value = 1304.2307
$ws_GSM.Range("I122").Value = 1345.5
$ws_GSM.Range("J122").Value = 1166.6666
$ws_GSM.Range("K122").Value = 4036.5
$ws_GSM.Range("L122").Value = 3499.9998
$ws_GSM.Range("M122").Value = -1586.5
$ws_GSM.Range("N122").Value = -8399.9998

# Row 48 (LTW) - hunk 23
$ws_LTW.Range("H48").Value = 13500
$ws_LTW.Range("I48").Value = 0
$ws_LTW.Range("J48").Value = 13500
$ws_LTW.Range("K48").Value = 0
$ws_LTW.Range("L48").Value = 13500
$ws_LTW.Range("M48").ClearContents()
$ws_LTW.Range("N48").Value = -14822

# Row 136 (LTW) - hunk 24
$ws_LTW.Range("H136").Value = 4207.485
$ws_LTW.Range("I136").Value = 2687.35
$ws_LTW.Range("J136").Value = 6546.154
$ws_LTW.Range("K136").Value = 8062.049999999999
$ws_LTW.Range("L136").Value = 19638.462
$ws_LTW.Range("M136").Value = -5512.049999999999
$ws_LTW.Range("N136").Value = -24738.462

# Row 122 (WVR) - hunk 25
$ws_WVR.Range("H122").Value = 30304144
$ws_WVR.Range("I122").Value = 34483864
$ws_WVR.Range("J122").Value = 1163.75
$ws_WVR.Range("K122").Value = 103451592
$ws_WVR.Range("L122").Value = 3491.25
$ws_WVR.Range("M122").Value = -103449142
$ws_WVR.Range("N122").Value = -8391.25

# Row 136 (WVR) - hunk 26
$ws_WVR.Range("H136").Value = 1313.4482
$ws_WVR.Range("I136").Value = 1471.25
$ws_WVR.Range("J136").Value = 556
$ws_WVR.Range("K136").Value = 4413.75
$ws_WVR.Range("L136").Value = 1668
$ws_WVR.Range("M136").Value = -1863.75
$ws_WVR.Range("N136").Value = -6768
